$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.318.81'
$ws.Range('E2').Value = '  -0.64%  '
$ws.Range('D3').Value = '1.713.84'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5291'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.006'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('E8').Value = '  +1.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2647'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.84'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07713'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.487'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.56%  '
$ws.Range('D13').Value = '1.949.02'
$ws.Range('E13').Value = '  -0.66%  '
$ws.Range('D14').Value = '1.715.55'
$ws.Range('E14').Value = '  -0.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5790'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').Value = '0.0₅8179'
$ws.Range('E16').Value = '  -1.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.70'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').Value = '27.341.25'
$ws.Range('E18').Value = '  -0.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '219.88'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.007'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.649'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.53%  '
$ws.Range('E22').Value = '  -1.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.025'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.26'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.709'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1208'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.237'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05384'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.49%  '
$ws.Range('E31').Value = '  -0.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.477'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.392'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.634'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.846'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9519'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.400'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5885'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.62%  '
$ws.Range('D39').Value = '1.154.91'
$ws.Range('E39').Value = '  +10.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01652'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.838'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.19%  '
$ws.Range('E42').Value = '  +0.28%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8392'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.89'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.34%  '
$ws.Range('D45').Value = '1.856.25'
$ws.Range('E45').Value = '  -0.66%  '
$ws.Range('D46').Value = '0.0₈118'
$ws.Range('E46').Value = '  +1.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '57.68'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.63%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4561'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.29%  '
$ws.Range('B49').Value = 'Frax'
$ws.Range('C49').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.004'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.110'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('E51').Value = '  -1.10%  '
